$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.226.66"
$ws.Range("E2").Value = "  +1.11%  "

$ws.Range("D3").Value = "3.688.84"
$ws.Range("E3").Value = "  +0.57%  "

$s_D4 = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $s_D4
$ws.Range("E4").Value = "  +0.16%  "

$s_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'606.78"
$ws.Range("D5").Style = $s_D5
$ws.Range("E5").Value = "  +4.76%  "

$s_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'192.73"
$ws.Range("D6").Style = $s_D6
$ws.Range("E6").Value = "  +13.60%  "

$s_D7 = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.630"
$ws.Range("D7").Style = $s_D7
$ws.Range("E7").Value = "  +1.68%  "

$ws.Range("E8").Value = "  +0.00%  "

$s_D9 = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.719"
$ws.Range("D9").Style = $s_D9
$ws.Range("E9").Value = "  +3.11%  "

$s_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'59.47"
$ws.Range("D10").Style = $s_D10
$ws.Range("E10").Value = "  +17.44%  "

$s_D11 = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.159"
$ws.Range("D11").Style = $s_D11
$ws.Range("E11").Value = "  -0.65%  "

$s_D12 = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0000284"
$ws.Range("D12").Style = $s_D12
$ws.Range("E12").Value = "  -0.38%  "

$s_D13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'10.38"
$ws.Range("D13").Style = $s_D13
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "4.265.81"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").Value = "3.679.96"
$ws.Range("E15").Value = "  +0.41%  "

$ws.Range("E16").Value = "  +1.22%  "

$s_D17 = $ws.Range("D17").Style
$ws.Range("D17").Value = "'19.30"
$ws.Range("D17").Style = $s_D17
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("E18").Value = "  +2.88%  "

$s_D19 = $ws.Range("D19").Style
$ws.Range("D19").Value = "'12.77"
$ws.Range("D19").Style = $s_D19
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").Value = "68.001.86"
$ws.Range("E20").Value = "  +1.08%  "

$s_D21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'406.23"
$ws.Range("D21").Style = $s_D21
$ws.Range("E21").Value = "  +0.75%  "

$s_D22 = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.56"
$ws.Range("D22").Style = $s_D22
$ws.Range("E22").Value = "  +2.65%  "

$s_D23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'89.70"
$ws.Range("D23").Style = $s_D23
$ws.Range("E23").Value = "  +2.87%  "

$s_D24 = $ws.Range("D24").Style
$ws.Range("D24").Value = "'11.55"
$ws.Range("D24").Style = $s_D24
$ws.Range("E24").Value = "  +8.89%  "

$s_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'3.05"
$ws.Range("D25").Style = $s_D25
$ws.Range("E25").Value = "  +1.17%  "

$s_D26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'12.94"
$ws.Range("D26").Style = $s_D26
$ws.Range("E26").Value = "  +2.34%  "

$s_D27 = $ws.Range("D27").Style
$ws.Range("D27").Value = "'6.02"
$ws.Range("D27").Style = $s_D27
$ws.Range("E27").Value = "  +0.52%  "

$s_D28 = $ws.Range("D28").Style
$ws.Range("D28").Value = "'3.75"
$ws.Range("D28").Style = $s_D28
$ws.Range("E28").Value = "  +0.67%  "

$s_D29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'9.56"
$ws.Range("D29").Style = $s_D29
$ws.Range("E29").Value = "  +2.07%  "

$s_D30 = $ws.Range("D30").Style
$ws.Range("D30").Value = "'32.42"
$ws.Range("D30").Style = $s_D30
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("E31").Value = "  +2.10%  "

$s_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'47.23"
$ws.Range("D32").Style = $s_D32
$ws.Range("E32").Value = "  +10.03%  "

$s_D33 = $ws.Range("D33").Style
$ws.Range("D33").Value = "'12.66"
$ws.Range("D33").Style = $s_D33
$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("E34").Value = "  +5.36%  "

$s_D35 = $ws.Range("D35").Style
$ws.Range("D35").Value = "'631.37"
$ws.Range("D35").Style = $s_D35
$ws.Range("E35").Value = "  +7.38%  "

$s_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'67.52"
$ws.Range("D36").Style = $s_D36
$ws.Range("E36").Value = "  +5.13%  "

$ws.Range("D37").Value = "0.0₃0825"
$ws.Range("E37").Value = "  -6.47%  "

$s_D38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.409"
$ws.Range("D38").Style = $s_D38
$ws.Range("E38").Value = "  +4.76%  "

$ws.Range("E39").Value = "  +0.00%  "

$s_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = $s_D40
$ws.Range("E40").Value = "  +0.11%  "

$s_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.138"
$ws.Range("D41").Style = $s_D41
$ws.Range("E41").Value = "  +3.90%  "

$s_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'3.01"
$ws.Range("D42").Style = $s_D42
$ws.Range("E42").Value = "  +1.90%  "

$s_D43 = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.0441"
$ws.Range("D43").Style = $s_D43
$ws.Range("E43").Value = "  +2.50%  "

$s_D44 = $ws.Range("D44").Style
$ws.Range("D44").Value = "'2.62"
$ws.Range("D44").Style = $s_D44
$ws.Range("E44").Value = "  -2.15%  "

$ws.Range("D45").Value = "2.885.99"
$ws.Range("E45").Value = "  +4.72%  "

$ws.Range("E46").Value = "  +4.96%  "

$s_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'9.16"
$ws.Range("D47").Style = $s_D47
$ws.Range("E47").Value = "  +0.92%  "

$s_D48 = $ws.Range("D48").Style
$ws.Range("D48").Value = "'145.89"
$ws.Range("D48").Style = $s_D48
$ws.Range("E48").Value = "  +3.37%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$s_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'2.67"
$ws.Range("D49").Style = $s_D49
$ws.Range("E49").Value = "  -5.16%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$s_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'2.65"
$ws.Range("D50").Style = $s_D50
$ws.Range("E50").Value = "  -0.03%  "

$s_D51 = $ws.Range("D51").Style
$ws.Range("D51").Value = "'3.05"
$ws.Range("D51").Style = $s_D51
$ws.Range("E51").Value = "  -2.31%  "
